$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create row 26 first, copying the formatting of row 25 (A column has a
# style applied, B column uses the default style) so the new cells end up
# styled the same way as the rest of the table.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Updated "pojazdy" (vehicle) counts per line (column A), and the shifted
# stop id list in column B after the extra line was added.
$ws.Range("A2").Value2  = 36
$ws.Range("B2").Value2  = 1

$ws.Range("A3").Value2  = 132
$ws.Range("B3").Value2  = 3

$ws.Range("A4").Value2  = 1
$ws.Range("B4").Value2  = 4

$ws.Range("A5").Value2  = 40
$ws.Range("B5").Value2  = 5

$ws.Range("A6").Value2  = 179
$ws.Range("B6").Value2  = 8

$ws.Range("A7").Value2  = 189
$ws.Range("B7").Value2  = 9

$ws.Range("A8").Value2  = 160
$ws.Range("B8").Value2  = 10

$ws.Range("A9").Value2  = 59
$ws.Range("B9").Value2  = 11

$ws.Range("A10").Value2 = 21
$ws.Range("B10").Value2 = 13

$ws.Range("A11").Value2 = 95
$ws.Range("B11").Value2 = 14

$ws.Range("A12").Value2 = 10
$ws.Range("B12").Value2 = 16

$ws.Range("A13").Value2 = 39
$ws.Range("B13").Value2 = 17

$ws.Range("A14").Value2 = 79
$ws.Range("B14").Value2 = 18

$ws.Range("A15").Value2 = 70
$ws.Range("B15").Value2 = 19

$ws.Range("A16").Value2 = 71
$ws.Range("B16").Value2 = 20

$ws.Range("A17").Value2 = 0
$ws.Range("B17").Value2 = 21

$ws.Range("A18").Value2 = 163
$ws.Range("B18").Value2 = 22

$ws.Range("A19").Value2 = 184
$ws.Range("B19").Value2 = 24

$ws.Range("A20").Value2 = 3
$ws.Range("B20").Value2 = 44

$ws.Range("A21").Value2 = 12
$ws.Range("B21").Value2 = 49

$ws.Range("A22").Value2 = 154
$ws.Range("B22").Value2 = 50

$ws.Range("A23").Value2 = 20
$ws.Range("B23").Value2 = 52

$ws.Range("A24").Value2 = 97
$ws.Range("B24").Value2 = 62

$ws.Range("A25").Value2 = 94
$ws.Range("B25").Value2 = 64

$ws.Range("A26").Value2 = 44
$ws.Range("B26").Value2 = 72
